$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (uppercased labels) ----
$ws.Cells.Item(1, 1).Value = "ITEM"
$ws.Cells.Item(1, 2).Value = "PREGUNTA"
$ws.Cells.Item(1, 3).Value = "ESCALA"
$ws.Cells.Item(1, 4).Value = "POSIBLES_RESPUESTAS"

# ---- Scale-text constants ----
$likert5 = "1: Totalmente en desacuerdo, 2: En desacuerdo, 3: Neutral, 4: De acuerdo, 5: Totalmente de acuerdo"
$bin2    = " 1: De acuerdo, 2: Totalmente de acuerdo"
$tri3    = "1: SI, 2: NO, 3: NO SE"

# For every data row (2-21): column C becomes a numeric count of possible
# answers (was text "likert" / "Binario"), and column D is updated so its
# text matches that numeric scale.
$escala = @{
    2  = 5
    3  = 2
    4  = 5
    5  = 2
    6  = 5
    7  = 2
    8  = 5
    9  = 5
    10 = 3
    11 = 5
    12 = 5
    13 = 5
    14 = 5
    15 = 5
    16 = 5
    17 = 5
    18 = 5
    19 = 5
    20 = 5
    21 = 3
}

# Rows whose "posibles respuestas" cell (column D) should be left-aligned
# (style already present on D2/D3; D5/D7 newly pick it up here).
$leftAlignRows = @(2, 3, 5, 7)

foreach ($r in 2..21) {
    $n = $escala[$r]
    $ws.Cells.Item($r, 3).Value = $n

    if ($n -eq 5) {
        $text = $likert5
    } elseif ($n -eq 2) {
        $text = $bin2
    } else {
        $text = $tri3
    }
    $ws.Cells.Item($r, 4).Value = $text

    if ($leftAlignRows -contains $r) {
        $ws.Cells.Item($r, 4).HorizontalAlignment = -4131
    }
}

# ---- Column B width widened ----
$ws.Range("B1").EntireColumn.ColumnWidth = 93.3

# ---- Selection / scroll position ----
$ws.Range("D2").Select()
$excel.ActiveWindow.ScrollColumn = 3
